# The sheet currently has yearly rows for 2007-2018 (data rows) followed by
# 2019 and 2020 (rows with only a total in column J). We need to:
#   1. Remove the oldest three years (2007, 2008, 2009) which occupy rows 2-4.
#      This shifts everything up by three rows, so 2010 becomes row 2 and the
#      last existing row (2020) becomes row 12.
#   2. Append a new row for 2021, following the same "only a J total" pattern
#      used by the 2019/2020 rows, with a total of 2932.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Delete the rows for 2007, 2008 and 2009.
$ws.Range("A2:J4").EntireRow.Delete()

# 2) Add the new trailing row for 2021 (now row 13), matching the style of
#    the row above it (2020, row 12) and leaving the detail columns blank.
$ws.Range("A12").Copy()
$ws.Range("A13").PasteSpecial(-4122)
$ws.Range("A13").Value = "2021年"
$ws.Range("J13").Value = 2932
